$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the target range to Text format first so numeric-looking strings
# (e.g. "220.66", "0.620", "1.20") are stored verbatim as text instead of
# being auto-coerced into numbers (which would also drop trailing zeros).
$rng = $ws.Range("D2:E51")
$rng.NumberFormat = "@"

$ws.Range('D2').Value = '33.812.55'
$ws.Range('E2').Value = '  -2.37%  '
$ws.Range('D3').Value = '1.767.82'
$ws.Range('E3').Value = '  -1.34%  '
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').Value = '220.66'
$ws.Range('E5').Value = '  -1.88%  '
$ws.Range('D6').Value = '0.548'
$ws.Range('E6').Value = '  -0.77%  '
$ws.Range('E8').Value = '  -5.22%  '
$ws.Range('E9').Value = '  -0.14%  '
$ws.Range('D10').Value = '0.0704'
$ws.Range('E10').Value = '  +5.66%  '
$ws.Range('E11').Value = '  -1.67%  '
$ws.Range('D12').Value = '2.023.39'
$ws.Range('D13').Value = '1.774.83'
$ws.Range('E13').Value = '  -0.75%  '
$ws.Range('D14').Value = '10.45'
$ws.Range('E14').Value = '  -5.58%  '
$ws.Range('D15').Value = '0.620'
$ws.Range('E15').Value = '  -2.02%  '
$ws.Range('D16').Value = '33.861.12'
$ws.Range('E16').Value = '  -2.29%  '
$ws.Range('E17').Value = '  -2.21%  '
$ws.Range('D18').Value = '67.49'
$ws.Range('E18').Value = '  -2.15%  '
$ws.Range('D19').Value = '242.77'
$ws.Range('E19').Value = '  -4.28%  '
$ws.Range('D20').Value = '0.0₃0770'
$ws.Range('E20').Value = '  +1.32%  '
$ws.Range('E21').Value = '  -0.09%  '
$ws.Range('D22').Value = '10.47'
$ws.Range('E22').Value = '  +0.97%  '
$ws.Range('D23').Value = '4.03'
$ws.Range('E23').Value = '  -4.51%  '
$ws.Range('D24').Value = '2.08'
$ws.Range('E24').Value = '  -2.26%  '
$ws.Range('D25').Value = '157.05'
$ws.Range('E25').Value = '  -0.94%  '
$ws.Range('D26').Value = '16.27'
$ws.Range('E26').Value = '  -0.66%  '
$ws.Range('D27').Value = '6.93'
$ws.Range('E27').Value = '  -1.89%  '
$ws.Range('E28').Value = '  -2.27%  '
$ws.Range('E29').Value = '  -0.26%  '
$ws.Range('D30').Value = '0.0519'
$ws.Range('E30').Value = '  +0.55%  '
$ws.Range('E31').Value = '  -2.01%  '
$ws.Range('D32').Value = '1.20'
$ws.Range('E32').Value = '  +0.17%  '
$ws.Range('D33').Value = '3.47'
$ws.Range('E33').Value = '  -3.05%  '
$ws.Range('E34').Value = '  -3.10%  '
$ws.Range('D35').Value = '1.389.29'
$ws.Range('E35').Value = '  -3.84%  '
$ws.Range('E36').Value = '  -0.88%  '
$ws.Range('D37').Value = '0.630'
$ws.Range('E37').Value = '  +0.82%  '
$ws.Range('E38').Value = '  -2.05%  '
$ws.Range('D39').Value = '0.921'
$ws.Range('E39').Value = '  +2.55%  '
$ws.Range('D40').Value = '2.35'
$ws.Range('E40').Value = '  -0.47%  '
$ws.Range('D41').Value = '78.40'
$ws.Range('E41').Value = '  -5.26%  '
$ws.Range('E42').Value = '  -5.39%  '
$ws.Range('E43').Value = '  +1.27%  '
$ws.Range('D44').Value = '5.86'
$ws.Range('E44').Value = '  -1.38%  '
$ws.Range('D45').Value = '1.03'
$ws.Range('E45').Value = '  -1.53%  '
$ws.Range('E46').Value = '  -4.02%  '
$ws.Range('D47').Value = '1.919.59'
$ws.Range('E47').Value = '  -2.13%  '
$ws.Range('D48').Value = '103.52'
$ws.Range('E48').Value = '  -0.04%  '
$ws.Range('E49').Value = '  -0.64%  '
$ws.Range('D50').Value = '11.68'
$ws.Range('E50').Value = '  -1.56%  '
$ws.Range('D51').Value = '0.0₆0119'
$ws.Range('E51').Value = '  -3.59%  '

# Restore the original (default) cell formatting so no stray number-format
# style is left behind on the edited cells.
$rng.ClearFormats()
